$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current last row (row 61), shifting the
# "Accuracy over PyType" row down to row 62.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the Scalpel accuracy summary.
$ws.Cells.Item(61, 3).Value = "Scalpel Accuracy:"
$ws.Cells.Item(61, 4).Value = 222.22
